$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 text with new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.33 = 13019.97 pesos`n✅ 13019.97 pesos = 3.33 = 965.8 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update "tasas" sheet N10, O10, N12, O12 values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 300
$wsTasas.Range("O10").Value = 3905.99
$wsTasas.Range("N12").Value = 3909.5
$wsTasas.Range("O12").Value = 290
